$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45923.01041666666, 0),
    @(3, 45923.02083333334, 0),
    @(4, 45923.03125, 0),
    @(5, 45923.04166666666, 0),
    @(6, 45923.05208333334, 0),
    @(7, 45923.0625, 0),
    @(8, 45923.07291666666, 0),
    @(9, 45923.08333333334, 0),
    @(10, 45923.09375, 0),
    @(11, 45923.10416666666, 0),
    @(12, 45923.11458333334, 0),
    @(13, 45923.125, 0),
    @(14, 45923.13541666666, 0),
    @(15, 45923.14583333334, 0),
    @(16, 45923.15625, 0),
    @(17, 45923.16666666666, 0),
    @(18, 45923.17708333334, 1),
    @(19, 45923.1875, 1),
    @(20, 45923.19791666666, 1),
    @(21, 45923.20833333334, 1),
    @(22, 45923.21875, 3),
    @(23, 45923.22916666666, 3),
    @(24, 45923.23958333334, 3),
    @(25, 45923.25, 3),
    @(26, 45923.26041666666, 64),
    @(27, 45923.27083333334, 72),
    @(28, 45923.28125, 88),
    @(29, 45923.29166666666, 104),
    @(30, 45923.30208333334, 511),
    @(31, 45923.3125, 549),
    @(32, 45923.32291666666, 593),
    @(33, 45923.33333333334, 644),
    @(34, 45923.34375, 1234),
    @(35, 45923.35416666666, 1283),
    @(36, 45923.36458333334, 1336),
    @(37, 45923.375, 1390),
    @(38, 45923.38541666666, 1843),
    @(39, 45923.39583333334, 1883),
    @(40, 45923.40625, 1924),
    @(41, 45923.41666666666, 1961),
    @(42, 45923.42708333334, 2212),
    @(43, 45923.4375, 2236),
    @(44, 45923.44791666666, 2260),
    @(45, 45923.45833333334, 2278),
    @(46, 45923.46875, 2377),
    @(47, 45923.47916666666, 2387),
    @(48, 45923.48958333334, 2396),
    @(49, 45923.5, 2400),
    @(50, 45923.51041666666, 2405),
    @(51, 45923.52083333334, 2404),
    @(52, 45923.53125, 2402),
    @(53, 45923.54166666666, 2397),
    @(54, 45923.55208333334, 2326),
    @(55, 45923.5625, 2314),
    @(56, 45923.57291666666, 2300),
    @(57, 45923.58333333334, 2285),
    @(58, 45923.59375, 2104),
    @(59, 45923.60416666666, 2079),
    @(60, 45923.61458333334, 2048),
    @(61, 45923.625, 2014),
    @(62, 45923.63541666666, 1655),
    @(63, 45923.64583333334, 1609),
    @(64, 45923.65625, 1566),
    @(65, 45923.66666666666, 1520),
    @(66, 45923.67708333334, 928),
    @(67, 45923.6875, 869),
    @(68, 45923.69791666666, 820),
    @(69, 45923.70833333334, 769),
    @(70, 45923.71875, 273),
    @(71, 45923.72916666666, 239),
    @(72, 45923.73958333334, 208),
    @(73, 45923.75, 187),
    @(74, 45923.76041666666, 20),
    @(75, 45923.77083333334, 18),
    @(76, 45923.78125, 17),
    @(77, 45923.79166666666, 17),
    @(78, 45923.80208333334, 10),
    @(79, 45923.8125, 10),
    @(80, 45923.82291666666, 10),
    @(81, 45923.83333333334, 9),
    @(82, 45923.84375, 1),
    @(83, 45923.85416666666, 1),
    @(84, 45923.86458333334, 1),
    @(85, 45923.875, 1),
    @(86, 45923.88541666666, 1),
    @(87, 45923.89583333334, 1),
    @(88, 45923.90625, 1),
    @(89, 45923.91666666666, 1),
    @(90, 45923.92708333334, 0),
    @(91, 45923.9375, 0),
    @(92, 45923.94791666666, 0),
    @(93, 45923.95833333334, 0),
    @(94, 45923.96875, 0),
    @(95, 45923.97916666666, 0),
    @(96, 45923.98958333334, 0),
    @(97, 45924.0, 0)
)

foreach ($item in $data) {
    $row = $item[0]
    $aVal = $item[1]
    $bVal = $item[2]
    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal
}
